$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.301.27"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.89%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.679.63"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.86%  "

# Row 4
$ws.Range("E4").Value = "  +0.23%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "218.30"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.94%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5251"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.00%  "

# Row 7
$ws.Range("E7").Value = "  +0.21%  "

# Row 8
$ws.Range("E8").Value = "  +2.31%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06467"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.73%  "

# Row 10
$ws.Range("E10").Value = "  +2.17%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07512"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.28%  "

# Row 12
$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.687.72"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.30%  "

# Row 13
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.529"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.52%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5807"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.27%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.000008517"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.25%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.81"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.02%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.327.94"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.75%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "4.926"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.31%  "

# Row 19
$ws.Range("E19").Value = "  +0.24%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.88"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.79%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "189.96"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.42%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.205"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.02%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.008"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.19%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "145.40"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.12%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "7.826"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.83%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1254"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +4.40%  "

# Row 27
$ws.Range("E27").Value = "  +1.42%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.06453"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.29%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.359"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +4.79%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.324"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.68%  "

# Row 31
$ws.Range("B31").Value = "InternetComputer(DFINITY)"
$ws.Range("C31").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.607"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.95%  "

# Row 32
$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.596"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.17%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.668"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.31%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.030"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.65%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.6237"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.54%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.409"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.98%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.737"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.35%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.439"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +4.45%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01626"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.17%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.107.99"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.00%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8761"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.94%  "

# Row 42
$ws.Range("E42").Value = "  +0.66%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "100.70"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.25%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.831.44"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.98%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00000000109"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -4.72%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "56.98"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.55%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.186"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.29%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.005"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.13%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.05267"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.24%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.4296"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.19%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.087"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.67%  "
